# Fruta / hortaliza, semanal
# Insert a new weekly record at row 64 (Vega Monumental Concepcion - Mango),
# pushing the existing rows 64:85 down to 65:86, then populate the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 64; Excel shifts rows 64:85 -> 65:86
$ws.Rows(64).Insert()

# Populate the newly inserted row 64 with this week's data
$ws.Range("A64").Value = 11
$ws.Range("B64").Value = "Vega Monumental Concepción"
$ws.Range("C64").Value = "Bíobío"
$ws.Range("D64").Value = 44559
$ws.Range("E64").Value = 8
$ws.Range("F64").Value = "Fruta"
$ws.Range("G64").Value = 100108
$ws.Range("H64").Value = "Tropicales y subtropicales"
$ws.Range("I64").Value = 100108002
$ws.Range("J64").Value = "Mango"
$ws.Range("K64").Value = "Sin especificar"
$ws.Range("L64").Value = "Primera"
$ws.Range("M64").Value = 300
$ws.Range("N64").Value = 6000
$ws.Range("O64").Value = 6500
$ws.Range("P64").Value = 6333
$ws.Range("Q64").Value = "$/bandeja 4 kilos"
$ws.Range("R64").Value = "Perú"
$ws.Range("S64").Value = 1583
$ws.Range("T64").Value = 4
